$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.250.94'
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").Value = '''1.839.99'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("D4").Value = '''0.9991'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''240.48'
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("D6").Value = '''0.6260'
$ws.Range("E6").Value = '  -0.35%  '

$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").Value = '''0.07446'
$ws.Range("E8").Value = '  -2.82%  '

$ws.Range("D9").Value = '''0.2894'
$ws.Range("E9").Value = '  -0.94%  '

$ws.Range("D10").Value = '''24.32'
$ws.Range("E10").Value = '  -1.73%  '

$ws.Range("D11").Value = '''0.07718'

$ws.Range("D12").Value = '''1.839.01'
$ws.Range("E12").Value = '  -1.82%  '

$ws.Range("D13").Value = '''4.987'

$ws.Range("D14").Value = '''0.6789'
$ws.Range("E14").Value = '  -0.42%  '

$ws.Range("D15").Value = '''0.00001015'
$ws.Range("E15").Value = '  -3.01%  '

$ws.Range("D16").Value = '''82.07'
$ws.Range("E16").Value = '  -1.74%  '

$ws.Range("D17").Value = '''2.097.77'
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("D18").Value = '''6.134'
$ws.Range("E18").Value = '  -0.87%  '

$ws.Range("D19").Value = '''29.281.18'
$ws.Range("E19").Value = '  -0.83%  '

$ws.Range("D20").Value = '''228.57'
$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("D21").Value = '''12.28'
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").Value = '''0.9999'
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").Value = '''7.371'
$ws.Range("E23").Value = '  -1.42%  '

$ws.Range("D24").Value = '''1.000'
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("D25").Value = '''158.47'
$ws.Range("E25").Value = '  +0.71%  '

$ws.Range("D26").Value = '''0.1377'
$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("D27").Value = '''8.378'
$ws.Range("E27").Value = '  -0.52%  '

$ws.Range("D28").Value = '''17.54'
$ws.Range("E28").Value = '  -1.29%  '

$ws.Range("D29").Value = '''1.398'
$ws.Range("E29").Value = '  +2.54%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.06136'
$ws.Range("E30").Value = '  +9.22%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.475'
$ws.Range("E31").Value = '  +0.60%  '

$ws.Range("D32").Value = '''4.091'
$ws.Range("E32").Value = '  -1.02%  '

$ws.Range("D33").Value = '''4.047'
$ws.Range("E33").Value = '  -0.10%  '

$ws.Range("D34").Value = '''1.821'
$ws.Range("E34").Value = '  -1.21%  '

$ws.Range("E35").Value = '  -1.93%  '

$ws.Range("D36").Value = '''0.7085'
$ws.Range("E36").Value = '  +0.94%  '

$ws.Range("D37").Value = '''2.586'
$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("D38").Value = '''2.822'
$ws.Range("E38").Value = '  +2.65%  '

$ws.Range("D39").Value = '''1.244.07'
$ws.Range("E39").Value = '  +1.38%  '

$ws.Range("D40").Value = '''0.01812'
$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("D41").Value = '''6.531'
$ws.Range("E41").Value = '  +1.32%  '

$ws.Range("D42").Value = '''0.9118'
$ws.Range("E42").Value = '  +0.71%  '

$ws.Range("D43").Value = '''0.9986'
$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("D44").Value = '''1.998.31'
$ws.Range("E44").Value = '  -0.82%  '

$ws.Range("D45").Value = '''101.68'
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("D46").Value = '''66.12'
$ws.Range("E46").Value = '  -0.12%  '

$ws.Range("D47").Value = '''7.041'
$ws.Range("E47").Value = '  -1.77%  '

$ws.Range("D48").Value = '''0.1161'
$ws.Range("E48").Value = '  +0.41%  '

$ws.Range("D49").Value = '''9.030'
$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '''0.00000000115'
$ws.Range("E50").Value = '  -3.02%  '

$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '''0.3941'
$ws.Range("E51").Value = '  -1.99%  '
